# Replace the sample student roster with a new set of rows.
# Column layout: A=name, C=admission no, E=email (hyperlinked), G/H=flags.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; Name = "Anex";      AdmNo = 1522; Email = "anex@gmail.com" },
    @{ Row = 3; Name = "Jeffin";    AdmNo = 9638; Email = "jeffin@gmail.com" },
    @{ Row = 4; Name = "Mahadevan"; AdmNo = 7531; Email = "mahadevan@gmail.com" },
    @{ Row = 5; Name = "Melvin";    AdmNo = 5548; Email = "Melvin@gmail.com" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Name
    $ws.Range("C$n").Value = $r.AdmNo
    # The cell already carries a mailto hyperlink; only the displayed text
    # changes, the underlying hyperlink address is left as-is.
    $ws.Range("E$n").Value = $r.Email
    $ws.Range("G$n").Value = 1
    $ws.Range("H$n").Value = 2
}

# Move the active selection to E4:F4 (matches the latest edited row).
$ws.Range("E4:F4").Select() | Out-Null
